$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if ($rng.Find.Found) {
        # Pin the run boundaries on both sides of the match with
        # temporary bookmarks before mutating the text. This stops the
        # engine from silently coalescing the edited run into an
        # adjacent run that happens to share identical formatting,
        # which would otherwise merge two distinct <w:r> elements into
        # one on save.
        $boundaryStart = $d.Range($rng.Start, $rng.Start)
        $boundaryEnd = $d.Range($rng.End, $rng.End)
        $d.Bookmarks.Add("zzBoundaryStart", $boundaryStart)
        $d.Bookmarks.Add("zzBoundaryEnd", $boundaryEnd)

        $rng.Text = $newText

        $d.Bookmarks("zzBoundaryStart").Delete()
        $d.Bookmarks("zzBoundaryEnd").Delete()
    }
}

Replace-ExactText "Passeport" "Carte d'identité nationale"
Replace-ExactText "  N°PP25342A  " "  N°AA-45467776-AQ  "
Replace-ExactText "14 mars 2019" "12 juillet 2023"
Replace-ExactText "Direction générale de la documentation et l'immigation" "Forces nationales de police"
Replace-ExactText "04 décembre 2024" "11 décembre 2024"
